# Handle separate tenant names for Cloud and On-Premises Orchestrator instances
#
# This script edits the "Settings" sheet (Table14) of the Config workbook so that
# the single "TenantName" row is split into a "CloudTenantName" row and an
# "OnPremisesTenantName" row, a new "CloudAccountName" row is (re)inserted right
# after the Cloud tenant row, the on-premises Orchestrator version sample value
# changes, and the Cloud Platform URL domain changes from platform.uipath.com to
# cloud.uipath.com (on the "Advanced Settings" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Step 1: compact the sheet by removing the two blank spacer rows -------
# Original layout (1-indexed rows):
#   1 Name | Value | Explanation                 (header)
#   2 EntitiesWorkbooksFolderPath
#   3 TenantName
#   4 (blank)
#   5 OnPremisesOrchestratorURL
#   6 OnPremisesOrchestratorVersion
#   7 (blank)
#   8 CloudAccountName
# Delete from the bottom up so row indexes above the deleted row stay valid.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(4).Delete()

# Sheet is now contiguous:
#   1 Name | Value | Explanation
#   2 EntitiesWorkbooksFolderPath
#   3 TenantName
#   4 OnPremisesOrchestratorURL
#   5 OnPremisesOrchestratorVersion
#   6 CloudAccountName

# --- Step 2: make room for the new CloudAccountName row under the (soon to
# be renamed) Cloud tenant row -----------------------------------------------
$ws.Rows.Item(4).Insert()

# Sheet now has 7 rows:
#   1 Name | Value | Explanation
#   2 EntitiesWorkbooksFolderPath
#   3 TenantName                    -> becomes CloudTenantName
#   4 (new, blank)                  -> becomes CloudAccountName
#   5 OnPremisesOrchestratorURL
#   6 OnPremisesOrchestratorVersion
#   7 CloudAccountName              -> becomes OnPremisesTenantName

# --- Step 3: row heights -----------------------------------------------------
$ws.Rows.Item(3).RowHeight = 58
$ws.Rows.Item(4).RowHeight = 72.5
$ws.Rows.Item(5).RowHeight = 58
$ws.Rows.Item(6).RowHeight = 58
$ws.Rows.Item(7).RowHeight = 159.5

# --- Step 4: cell values ------------------------------------------------------

# Row 3: CloudTenantName
$ws.Cells.Item(3, 1).Value = "CloudTenantName"
$ws.Cells.Item(3, 2).Value = "Default"
$ws.Cells.Item(3, 3).Value = "Name of the tenant to be used in case of Automation Cloud Orchestrator instances. " + [char]10 + "Sample value: Default."

# Row 4: CloudAccountName (new row)
$ws.Cells.Item(4, 1).Value = "CloudAccountName"
$ws.Cells.Item(4, 2).Value = "SampleAccount"
$ws.Cells.Item(4, 2).Style = $ws.Cells.Item(3, 2).Style
$ws.Cells.Item(4, 3).Value = "Unique site URL for Automation Cloud organization." + [char]10 + "This parameter is exclusive to Automation Cloud Orchestrator instances." + [char]10 + "Sample value: SampleAccount (assuming the organization account URL is https://cloud.uipath.com/SampleAccount)."

# Row 5: OnPremisesOrchestratorURL (unchanged content, just shifted)
$ws.Cells.Item(5, 1).Value = "OnPremisesOrchestratorURL"
$ws.Cells.Item(5, 2).Value = "https://myOrchestratorURL"
$ws.Cells.Item(5, 3).Value = "URL of the Orchestrator instance to be used. " + [char]10 + "This parameter is exclusive to on-premises Orchestrator instances." + [char]10 + [char]10 + "Sample value: https://myOrchestratorURL"

# Row 6: OnPremisesOrchestratorVersion (sample value number changes)
$ws.Cells.Item(6, 1).Value = "OnPremisesOrchestratorVersion"
$ws.Cells.Item(6, 2).Value = 201804
$ws.Cells.Item(6, 3).Value = "Version of the Orchestrator instance to be used. " + [char]10 + "This parameter is exclusive to on-premises Orchestrator instances." + [char]10 + "The specified value must be in the form YYYYMM, where YYYY is the 4-digit representation of an year and MM is the 2-digit representation of a month." + [char]10 + [char]10 + "Supported values: " + [char]10 + "201804" + [char]10 + "201904" + [char]10 + "201910" + [char]10 + "202004"

# Row 7: OnPremisesTenantName
$ws.Cells.Item(7, 1).Value = "OnPremisesTenantName"
$ws.Cells.Item(7, 2).Value = "Default"
$ws.Cells.Item(7, 3).Value = "Name of the tenant to be used in case of on-premises Orchestrator instances. " + [char]10 + "Sample value: Default."

# --- Step 5: fix up dimension / table range ----------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C7"))

# --- Step 6: Cloud Platform URL domain change (Advanced Settings sheet) -----
$ws2 = $wb.Worksheets.Item("Advanced Settings")
$found = $false
for ($r = 1; $r -le $ws2.UsedRange.Rows.Count; $r++) {
    if ($ws2.Cells.Item($r, 2).Value2 -eq "https://platform.uipath.com/") {
        $ws2.Cells.Item($r, 2).Value = "https://cloud.uipath.com/"
        $found = $true
        break
    }
}

# --- Step 7: tidy up the duplicated selection reference on the other sheets -
$ws2.Range("A2").Select()
$ws3 = $wb.Worksheets.Item("Localization")
$ws3.Range("A2").Select()
$ws.Select()

Write-Output "done; cloud url replaced: $found"
